$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Title paragraph: "Project Plan " -> "Project Plan Group 33"
#    (added as its own run, right after the existing trailing-space run)
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$r1 = $p1.Range
if ($r1.Text.TrimEnd([char]13) -eq "Project Plan ") {
    $insPoint1 = $d.Range($r1.End - 1, $r1.End - 1)
    $insPoint1.InsertAfter("Group 33")
}

# ---------------------------------------------------------------------------
# 2) Add a blank paragraph followed by an "Add the name?" paragraph right
#    after the "Student Names" paragraph (and before the page-break
#    paragraph that precedes the Table of Contents).
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$r3 = $p3.Range
if ($r3.Text.TrimEnd([char]13) -eq "Student Names") {
    $r3.InsertParagraphAfter()
    $p4 = $d.Paragraphs.Item(4)
    $p4.Range.InsertParagraphAfter()
    $p5 = $d.Paragraphs.Item(5)
    $p5.Range.InsertBefore("Add the name?")
}

# ---------------------------------------------------------------------------
# 3) Split the run ending in "...included in the WBS" so that "WBS" becomes
#    its own run (mirrors Word wrapping it with a gramStart/gramEnd proofing
#    span when the grammar checker flags it).
# ---------------------------------------------------------------------------
$wbsRange = $d.Content
$found = $wbsRange.Find.Execute("anything you do or will need to do should be included in the WBS", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $wbsOnly = $d.Range($wbsRange.End - 3, $wbsRange.End)
    if ($wbsOnly.Text -eq "WBS") {
        # Toggling a formatting property forces Word to split the run at the
        # sub-range boundary without altering the visible formatting.
        $wbsOnly.Font.Bold = $true
        $wbsOnly.Font.Bold = $false
    }
}

Write-Output "Edits applied"
